# Active_Outages.xlsx update - 6/18/2025, 5:00:57 PM
# - Refreshes "Elapsed Duration(Hrs)" values across R1, R2, R4, R5, R6
# - Updates the in-progress R4 outage (row 5 on R1) to a new Haj-removal
#   ticket now tracked under region R5 / power source SCECO / status "Good"
# - Adds that same new R5 outage record as a new row on the R3 sheet

$wb = $excel.ActiveWorkbook

# --- R1 sheet ---
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3930:15:01"
$ws1.Range("G3").Value = "69:47:39"
$ws1.Range("G4").Value = "92:47:39"
$ws1.Range("B5").Value = "R5"
$ws1.Range("D5").Value = "HAJ0155"
$ws1.Range("I5").Value = "SCECO"
$ws1.Range("J5").Value = "Good"

# --- R2 sheet ---
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12111:38:28"
$ws2.Range("G3").Value = "3241:21:57"
$ws2.Range("G4").Value = "479:33:31"

# --- R3 sheet: append new outage row (R5 / HAJ0125 / Haj Removal) ---
$ws3 = $wb.Worksheets.Item("R3")
$ws3.Range("B3").Value = "R5"
$ws3.Range("D3").Value = "HAJ0125"
$ws3.Range("I3").Value = "SCECO"
$ws3.Range("J3").Value = "Haj Removal"
$ws3.Range("L3").Value = "Latis"

# --- R4 sheet ---
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2957:28:17"
$ws4.Range("G3").Value = "184:40:32"
$ws4.Range("G4").Value = "72:52:57"
$ws4.Range("G5").Value = "70:30:30"

# --- R5 sheet ---
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "431:27:16"

# --- R6 sheet ---
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "71:59:34"
